# Sort each "Collection" block's rows (columns B:F) alphabetically by the
# Notation column (B), instead of by the PrefLabel column (C), as the
# original data had been sorted.  Column A (the merged Collection label)
# stays put since it belongs to the block as a whole, not to any one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count()
$lastCol = 6   # A..F

for ($r = 2; $r -le $lastRow; $r++) {
    $a = $ws.Cells.Item($r, 1).Value()
    if ($a -eq $null -or $a -eq "") {
        continue
    }

    # This row starts a merged "Collection" block; find its extent.
    $area = $ws.Cells.Item($r, 1).MergeArea
    $startRow = $area.Row
    $endRow = $startRow + $area.Rows.Count() - 1

    # Collect the B:F values for every row in the block.
    $blockRows = @()
    for ($rr = $startRow; $rr -le $endRow; $rr++) {
        $rowData = @{}
        for ($c = 2; $c -le $lastCol; $c++) {
            $rowData[$c] = $ws.Cells.Item($rr, $c).Value()
        }
        $blockRows += ,$rowData
    }

    # Sort by column B (Notation) value, ascending.
    $sortedRows = $blockRows | Sort-Object { $_[2] }

    # Write the sorted B:F values back into the block, row by row.
    for ($i = 0; $i -lt $sortedRows.Count; $i++) {
        $rr = $startRow + $i
        $rowData = $sortedRows[$i]
        for ($c = 2; $c -le $lastCol; $c++) {
            $ws.Cells.Item($rr, $c).Value = $rowData[$c]
        }
    }
}
